$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = "2024-07-19T08:00:00.000Z"
$ws.Range("D4:D32").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("D33:D73").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("D74:D110").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("D111:D127").Value = "2024-07-19T07:59:00.000Z"
